$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 values
$ws.Range("G9").Value = 5
$ws.Range("I9").Value = 7
$ws.Range("J9").Value = 4
$ws.Range("L9").Value = 3
$ws.Range("M9").Value = 3

# Update row 16 values
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 10
$ws.Range("I16").Value = 10
$ws.Range("J16").Value = 4
$ws.Range("K16").Value = 4
$ws.Range("L16").Value = 4
$ws.Range("M16").Value = 3

# Update row 17 values
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 10
$ws.Range("I17").Value = 9
$ws.Range("J17").Value = 4
$ws.Range("K17").Value = 4
$ws.Range("L17").Value = 4
$ws.Range("M17").Value = 3

# Update the frozen pane view / selection to match the scrolled position
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("H34").Select()
